# B6-PowerPoint.pptx edit
#
# 1) The three tables in the deck get their table style switched from
#    the custom "Table_0" style ({035C294F-DB96-45EE-8B7B-418866A2AEE8})
#    to the built-in "No Style, Table Grid" style
#    ({7044556C-4EDD-4A94-8C90-898FC3A7255E}).
#
# 2) The deck's applied design theme is switched from the "Integral"
#    ("Red Violet" colour scheme) back to the stock "Office Theme" -
#    every themed colour slot (dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink)
#    is reset to the default Office values.

$p = $ppt.ActivePresentation

# --- 1. Re-style every table in the deck -----------------------------
$oldStyle = "{035C294F-DB96-45EE-8B7B-418866A2AEE8}"
$newStyle = "{7044556C-4EDD-4A94-8C90-898FC3A7255E}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable -eq -1 -and $shape.Table.Style -eq $oldStyle) {
            $shape.Table.ApplyStyle($newStyle)
        }
    }
}

# --- 2. Restore the stock "Office" colour scheme on the design theme -
# RGB() packs as r + g*256 + b*65536, matching the VBA/PowerPoint
# ThemeColorScheme.Item(n).RGB convention.
function RGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$master = $p.Slides.Item(1).Master
$colors = $master.Theme.ThemeColorScheme

$colors.Item(1).RGB  = (RGB 0x00 0x00 0x00)   # dk1      000000
$colors.Item(2).RGB  = (RGB 0xFF 0xFF 0xFF)   # lt1      FFFFFF
$colors.Item(3).RGB  = (RGB 0x44 0x54 0x6A)   # dk2      44546A
$colors.Item(4).RGB  = (RGB 0xE7 0xE6 0xE6)   # lt2      E7E6E6
$colors.Item(5).RGB  = (RGB 0x5B 0x9B 0xD5)   # accent1  5B9BD5
$colors.Item(6).RGB  = (RGB 0xED 0x7D 0x31)   # accent2  ED7D31
$colors.Item(7).RGB  = (RGB 0xA5 0xA5 0xA5)   # accent3  A5A5A5
$colors.Item(8).RGB  = (RGB 0xFF 0xC0 0x00)   # accent4  FFC000
$colors.Item(9).RGB  = (RGB 0x44 0x72 0xC4)   # accent5  4472C4
$colors.Item(10).RGB = (RGB 0x70 0xAD 0x47)   # accent6  70AD47
$colors.Item(11).RGB = (RGB 0x05 0x63 0xC1)   # hlink    0563C1
$colors.Item(12).RGB = (RGB 0x95 0x4F 0x72)   # folHlink 954F72
